$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.121.21'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.52%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.601.88'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.85%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.33'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.83%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3780'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.77%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3642'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.15%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '50.15'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.66%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.255'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -6.67%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.003'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08131'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.71%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.01'
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.570'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -6.85%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001257'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.92%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.354'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -7.95%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.600.46'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.94%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.69'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.67%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06876'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.54%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.27'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -7.08%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.550'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.71%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("B22").Value = 'BitDAO'

$ws.Range("C22").Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.5570'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.67%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("B23").Value = 'Dai'

$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("B24").Value = 'Cosmos'

$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.99'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.44%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("B25").Value = 'WrappedBTC'

$ws.Range("C25").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '23.095.77'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.63%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").Value = 'Toncoin'

$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.346'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.38%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").Value = 'LidoDAOToken'

$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.728'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -7.15%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Value = 'EthereumClassic'

$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.06'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.45%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Value = 'Monero'

$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '149.96'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.10%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").Value = 'HuobiToken'

$ws.Range("C30").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.252'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.01%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = 'BitcoinCash'

$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '131.86'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.39%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = 'WEMIXTOKEN'

$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.429'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.43%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Value = 'Filecoin'

$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.812'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -12.63%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = 'WrappedliquidstakedEther2.0'

$ws.Range("C34").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.777.12'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.81%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").Value = 'ImmutableX'

$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9508'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -6.69%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = 'Hedera'

$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.07650'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.79%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = 'VeChain'

$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02732'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -6.54%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'

$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.223'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -7.89%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = 'Algorand'

$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2542'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.89%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = 'Stellar'

$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08906'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.71%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = 'FraxShare'

$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.03'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.20%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = 'TrustWalletToken'

$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.370'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.31%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = 'TheSandbox'

$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7073'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -6.64%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = 'Aptos'

$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.63'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.40%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = 'EnergySwap'

$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.46'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.27%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = 'Decentraland'

$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6600'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.15%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = 'Frax'

$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = 'NEARProtocol'

$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.294'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -6.19%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = 'PancakeSwap'

$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.980'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.73%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = 'Quant'

$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.62'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.01%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = 'Cronos'

$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07947'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.38%  '
$ws.Range("E51").Style = "Normal"
